# Updated cryptos list on Sun Feb 11 12:22:33 UTC 2024 with GitHub Actions
#
# Refreshes the "Price" (column D) and "Volume(1h)" (column E) figures for
# each coin row, and fixes the ranking swap between Solana/BNB (rows 5-6)
# and WEMIXToken/EnergySwap (rows 42-43). All cells in D/E hold plain text
# (not numbers), so for any new Price value that looks numeric we force the
# cell to Text format first - otherwise Excel's COM layer would silently
# re-interpret e.g. "266.00" or "2.00" as the number 266 / 2 and drop the
# original text formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '48.487.64'
$ws.Range('E2').Value = '  +2.58%  '
$ws.Range('D3').Value = '2.528.55'
$ws.Range('E3').Value = '  +1.59%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('B5').Value = 'Solana'
$ws.Range('C5').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '110.05'
$ws.Range('E5').Value = '  +2.02%  '
$ws.Range('B6').Value = 'BNB'
$ws.Range('C6').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '322.56'
$ws.Range('E6').Value = '  +0.61%  '
$ws.Range('E7').Value = '  +2.15%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('E9').Value = '  +3.45%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '40.59'
$ws.Range('E10').Value = '  +5.15%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.43'
$ws.Range('E11').Value = '  +12.23%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0822'
$ws.Range('E12').Value = '  +1.64%  '
$ws.Range('E13').Value = '  +1.11%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.29'
$ws.Range('E14').Value = '  +2.47%  '
$ws.Range('D15').Value = '2.923.95'
$ws.Range('E15').Value = '  +2.36%  '
$ws.Range('D16').Value = '2.534.71'
$ws.Range('E16').Value = '  +2.54%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.856'
$ws.Range('E17').Value = '  +1.27%  '
$ws.Range('D18').Value = '48.292.23'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.49'
$ws.Range('E19').Value = '  +6.04%  '
$ws.Range('E20').Value = '  +0.40%  '
$ws.Range('E21').Value = '  +1.90%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.70'
$ws.Range('E22').Value = '  -0.55%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '71.99'
$ws.Range('E23').Value = '  +2.33%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '266.00'
$ws.Range('E24').Value = '  +8.47%  '
$ws.Range('E25').Value = '  +1.09%  '
$ws.Range('E26').Value = '  -0.04%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '26.09'
$ws.Range('E27').Value = '  +1.54%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.37'
$ws.Range('E28').Value = '  +3.77%  '
$ws.Range('E29').Value = '  +1.49%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.145'
$ws.Range('E30').Value = '  +6.14%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '35.81'
$ws.Range('E31').Value = '  +2.98%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '49.72'
$ws.Range('E32').Value = '  +0.45%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '19.75'
$ws.Range('E33').Value = '  -1.59%  '
$ws.Range('E34').Value = '  +1.03%  '
$ws.Range('E35').Value = '  -0.03%  '
$ws.Range('E36').Value = '  +1.00%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.00'
$ws.Range('E37').Value = '  +1.93%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.71'
$ws.Range('E38').Value = '  +1.94%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.02'
$ws.Range('E39').Value = '  +2.97%  '
$ws.Range('E40').Value = '  +0.84%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '122.05'
$ws.Range('E41').Value = '  +2.24%  '
$ws.Range('B42').Value = 'WEMIXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.21'
$ws.Range('E42').Value = '  -0.04%  '
$ws.Range('B43').Value = 'EnergySwap'
$ws.Range('C43').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '21.95'
$ws.Range('E43').Value = '  +0.53%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0302'
$ws.Range('E44').Value = '  +2.47%  '
$ws.Range('D45').Value = '2.026.95'
$ws.Range('E45').Value = '  +2.13%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.17'
$ws.Range('E46').Value = '  +5.23%  '
$ws.Range('E47').Value = '  +8.43%  '
$ws.Range('E48').Value = '  +2.62%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.13'
$ws.Range('E49').Value = '  +0.73%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '5.24'
$ws.Range('E50').Value = '  +2.52%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '79.30'
$ws.Range('E51').Value = '  +3.04%  '
